$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column R (2021 data) that mirrors the formatting of the
# existing column Q for each row, then overwrite with the 2021 values.
$xlPasteFormats = -4122

$values = @{
    4  = 2021
    5  = 17.92
    6  = 19.65
    7  = 16.08
    8  = 16.86
    9  = 17.26
    10 = 16.44
    11 = 17.4
    12 = 18.1
    13 = 16.63
    14 = 26.64
    15 = 29.03
    16 = 24.15
    17 = 13.87
    18 = 13.87
    19 = 11.73
    20 = 11.98
    21 = 12.72
    22 = 11.2
    23 = 24.53
    24 = 31.24
    25 = 17.85
    26 = 13.54
    27 = 13.56
    28 = 13.52
    29 = 17.84
    30 = 20.85
    31 = 14.59
    32 = 30.69
    33 = 35.45
    34 = 25.64
}

foreach ($r in 4..34) {
    $srcCell = $ws.Range("Q$r")
    $dstCell = $ws.Range("R$r")

    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial($xlPasteFormats) | Out-Null

    $dstCell.Value = $values[$r]
}

$excel.CutCopyMode = 0

# Match the active cell/selection left behind in the authored workbook.
$ws.Range("S4").Select() | Out-Null
